$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (row 187), pushing the
# existing rows 187:210 down to 189:212. This matches the weekly refresh:
# two new observations (Brasil / Perú) for a newer date are prepended.
$ws.Rows("187:188").Insert()

# New row 187: Brasil, fecha 44474
$ws.Range("A187").Value = 3
$ws.Range("B187").Value = "Femacal de La Calera"
$ws.Range("C187").Value = "Coquimbo"
$ws.Range("D187").Value = 44474
$ws.Range("E187").Value = 5
$ws.Range("F187").Value = "Fruta"
$ws.Range("G187").Value = 100108
$ws.Range("H187").Value = "Tropicales y subtropicales"
$ws.Range("I187").Value = 100108002
$ws.Range("J187").Value = "Mango"
$ws.Range("K187").Value = "Sin especificar"
$ws.Range("L187").Value = "Primera"
$ws.Range("M187").Value = 228
$ws.Range("N187").Value = 9000
$ws.Range("O187").Value = 9000
$ws.Range("P187").Value = 9000
$ws.Range("Q187").Value = "$/bandeja 4 kilos"
$ws.Range("R187").Value = "Brasil"
$ws.Range("S187").Value = 2250
$ws.Range("T187").Value = 4

# New row 188: Perú, fecha 44474
$ws.Range("A188").Value = 3
$ws.Range("B188").Value = "Femacal de La Calera"
$ws.Range("C188").Value = "Coquimbo"
$ws.Range("D188").Value = 44474
$ws.Range("E188").Value = 5
$ws.Range("F188").Value = "Fruta"
$ws.Range("G188").Value = 100108
$ws.Range("H188").Value = "Tropicales y subtropicales"
$ws.Range("I188").Value = 100108002
$ws.Range("J188").Value = "Mango"
$ws.Range("K188").Value = "Sin especificar"
$ws.Range("L188").Value = "Primera"
$ws.Range("M188").Value = 228
$ws.Range("N188").Value = 9000
$ws.Range("O188").Value = 9000
$ws.Range("P188").Value = 9000
$ws.Range("Q188").Value = "$/bandeja 4 kilos"
$ws.Range("R188").Value = "Perú"
$ws.Range("S188").Value = 2250
$ws.Range("T188").Value = 4
